$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.467670917510986
$ws.Range("B1").Value = 1.654739737510681
$ws.Range("C1").Value = 2.02000880241394
$ws.Range("D1").Value = 2.725885629653931
$ws.Range("E1").Value = 6.623863220214844
